$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F12: new text, with wrap-text style (matches style used in E12/D12, s="4")
$ws.Range("F12").Value = "36 m²             5.40m²  linia"
$ws.Range("F12").WrapText = $true

# Update F14: new text
$ws.Range("F14").Value = "59.6 m²"

# Update active selection to I12 (matches recorded selection in diff)
$ws.Range("I12").Select()
